# Update stock rankings and weekly/monthly performance figures for the
# "Top Gainers" sheet of the market health workbook (data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

# Row 2
$ws.Range("C2").Value = 16.8578
$ws.Range("D2").Value = 16.2699
$ws.Range("E2").Value = 13.4666

# Row 3
$ws.Range("C3").Value = 14.2538
$ws.Range("D3").Value = 18.8437
$ws.Range("E3").Value = 19.6207

# Row 4
$ws.Range("C4").Value = 11.9847
$ws.Range("D4").Value = 8.0702
$ws.Range("E4").Value = 3.694

# Row 5
$ws.Range("C5").Value = 10.8463
$ws.Range("D5").Value = 10.8566
$ws.Range("E5").Value = 10.9392

# Row 8  # -> GENUSPOWER
$ws.Range("B8").Value = "GENUSPOWER"
$ws.Range("C8").Value = 9.2905
$ws.Range("D8").Value = 7.5287
$ws.Range("E8").Value = 4.338

# Row 9  # -> VBL
$ws.Range("B9").Value = "VBL"
$ws.Range("C9").Value = 9.138
$ws.Range("D9").Value = 7.4231
$ws.Range("E9").Value = 11.7084

# Row 10
$ws.Range("C10").Value = 8.938599999999999
$ws.Range("D10").Value = 16.1468
$ws.Range("E10").Value = 14.4404

# Row 11
$ws.Range("C11").Value = 7.9282
$ws.Range("D11").Value = 11.9255
$ws.Range("E11").Value = 14.3197

# Row 12
$ws.Range("C12").Value = 7.3441
$ws.Range("D12").Value = 7.8792
$ws.Range("E12").Value = 14.3149

# Row 13  # -> FISCHER
$ws.Range("B13").Value = "FISCHER"
$ws.Range("C13").Value = 7.3432
$ws.Range("D13").Value = 12.4028
$ws.Range("E13").Value = 5.4865

# Row 14  # -> BUTTERFLY
$ws.Range("B14").Value = "BUTTERFLY"
$ws.Range("C14").Value = 7.2525
$ws.Range("D14").Value = 10.1527
$ws.Range("E14").Value = 12.7428

# Row 15  # -> ABREL
$ws.Range("B15").Value = "ABREL"
$ws.Range("C15").Value = 7.0164
$ws.Range("D15").Value = 7.7924
$ws.Range("E15").Value = 7.3458

# Row 16
$ws.Range("C16").Value = 6.969
$ws.Range("D16").Value = 7.5968
$ws.Range("E16").Value = 16.5274

# Row 17
$ws.Range("C17").Value = 6.7308
$ws.Range("D17").Value = 3.1416
$ws.Range("E17").Value = 5.0947

# Row 18
$ws.Range("C18").Value = 6.5451
$ws.Range("D18").Value = 8.767200000000001
$ws.Range("E18").Value = 4.707

# Row 19
$ws.Range("C19").Value = 6.5011
$ws.Range("D19").Value = 4.0415
$ws.Range("E19").Value = -5.1468

# Row 20
$ws.Range("C20").Value = 6.4105
$ws.Range("D20").Value = -1.051
$ws.Range("E20").Value = 19.2608

# Row 21  # -> JISLJALEQS
$ws.Range("B21").Value = "JISLJALEQS"
$ws.Range("C21").Value = 6.1814
$ws.Range("D21").Value = 5.4736
$ws.Range("E21").Value = -0.628

# Row 22  # -> EPACKPEB
$ws.Range("B22").Value = "EPACKPEB"
$ws.Range("C22").Value = 6.1634
$ws.Range("D22").Value = -0.3984
$ws.Range("E22").Value = "N/A"

# Row 23
$ws.Range("C23").Value = 5.9839
$ws.Range("D23").Value = 9.7745
$ws.Range("E23").Value = 17.111

# Row 24  # -> GRAPHITE
$ws.Range("B24").Value = "GRAPHITE"
$ws.Range("C24").Value = 5.7683
$ws.Range("D24").Value = 12.2023
$ws.Range("E24").Value = 12.4043

# Row 25  # -> UTKARSHBNK
$ws.Range("B25").Value = "UTKARSHBNK"
$ws.Range("C25").Value = 5.6468
$ws.Range("D25").Value = -5.205
$ws.Range("E25").Value = -1.9066

# Row 26  # -> IOC
$ws.Range("B26").Value = "IOC"
$ws.Range("C26").Value = 5.6303
$ws.Range("D26").Value = 8.5456
$ws.Range("E26").Value = 8.9659

# Row 27  # -> PDSL
$ws.Range("B27").Value = "PDSL"
$ws.Range("C27").Value = 5.6095
$ws.Range("D27").Value = 3.6145
$ws.Range("E27").Value = 9.476699999999999

# Row 28
$ws.Range("C28").Value = 5.523
$ws.Range("D28").Value = 4.3602
$ws.Range("E28").Value = 14.2708

# Row 29  # -> SRM
$ws.Range("B29").Value = "SRM"
$ws.Range("C29").Value = 5.496
$ws.Range("D29").Value = 5.1825
$ws.Range("E29").Value = 6.1089

# Row 30  # -> ABDL
$ws.Range("B30").Value = "ABDL"
$ws.Range("C30").Value = 5.4666
$ws.Range("D30").Value = 4.3665
$ws.Range("E30").Value = 27.1125

# Row 31  # -> DATAMATICS
$ws.Range("B31").Value = "DATAMATICS"
$ws.Range("C31").Value = 5.2683
$ws.Range("D31").Value = 7.6915
$ws.Range("E31").Value = 16.1356

# Row 32  # -> STLTECH
$ws.Range("B32").Value = "STLTECH"
$ws.Range("C32").Value = 5.1955
$ws.Range("D32").Value = 2.1734
$ws.Range("E32").Value = 8.325200000000001

# Row 33
$ws.Range("C33").Value = 5.0098
$ws.Range("D33").Value = 2.4464
$ws.Range("E33").Value = 10.9028

# Row 39
$ws.Range("C39").Value = 4.7214
$ws.Range("D39").Value = 2.2659
$ws.Range("E39").Value = 30.4412

# Row 40  # -> GMBREW
$ws.Range("B40").Value = "GMBREW"
$ws.Range("C40").Value = 4.7013
$ws.Range("D40").Value = 0.2373
$ws.Range("E40").Value = 80.41

# Row 41  # -> LLOYDSENT
$ws.Range("B41").Value = "LLOYDSENT"
$ws.Range("C41").Value = 4.7002
$ws.Range("D41").Value = 1.966
$ws.Range("E41").Value = 11.3782

# Row 42  # -> FILATEX
$ws.Range("B42").Value = "FILATEX"
$ws.Range("C42").Value = 4.5595
$ws.Range("D42").Value = 9.948600000000001
$ws.Range("E42").Value = 25.631

# Row 43
$ws.Range("C43").Value = 4.5416
$ws.Range("D43").Value = 8.537100000000001
$ws.Range("E43").Value = 2.8017

# Row 44  # -> SURYAROSNI
$ws.Range("B44").Value = "SURYAROSNI"
$ws.Range("C44").Value = 4.4848
$ws.Range("D44").Value = 10.8573
$ws.Range("E44").Value = 2.5323

# Row 45
$ws.Range("C45").Value = 4.3502
$ws.Range("D45").Value = 3.328
$ws.Range("E45").Value = 3.2587

# Row 46  # -> HUDCO
$ws.Range("B46").Value = "HUDCO"
$ws.Range("C46").Value = 4.3068
$ws.Range("D46").Value = 3.8792
$ws.Range("E46").Value = 5.3749

# Row 47  # -> TCI
$ws.Range("B47").Value = "TCI"
$ws.Range("C47").Value = 4.2246
$ws.Range("D47").Value = 4.1281
$ws.Range("E47").Value = 4.6305

# Row 48  # -> GPPL
$ws.Range("B48").Value = "GPPL"
$ws.Range("C48").Value = 4.1952
$ws.Range("D48").Value = 3.1892
$ws.Range("E48").Value = 4.8282

# Row 49  # -> SUNFLAG
$ws.Range("B49").Value = "SUNFLAG"
$ws.Range("C49").Value = 4.1485
$ws.Range("D49").Value = 4.485
$ws.Range("E49").Value = 4.7837

# Row 50  # -> MRPL
$ws.Range("B50").Value = "MRPL"
$ws.Range("C50").Value = 4.1203
$ws.Range("D50").Value = 9.5589
$ws.Range("E50").Value = 19.8885

# Row 51  # -> APARINDS
$ws.Range("B51").Value = "APARINDS"
$ws.Range("C51").Value = 4.0791
$ws.Range("D51").Value = 8.536099999999999
$ws.Range("E51").Value = 15.7953

# Row 53  # -> WELSPUNLIV
$ws.Range("B53").Value = "WELSPUNLIV"
$ws.Range("C53").Value = 4.0514
$ws.Range("D53").Value = 4.1732
$ws.Range("E53").Value = 16.4342

# Row 54
$ws.Range("C54").Value = 4.0479
$ws.Range("D54").Value = 3.9777
$ws.Range("E54").Value = 3.2111

# Row 55  # -> SALASAR
$ws.Range("B55").Value = "SALASAR"
$ws.Range("C55").Value = 4.0042
$ws.Range("D55").Value = 5
$ws.Range("E55").Value = 11.274

# Row 56  # -> NBCC
$ws.Range("B56").Value = "NBCC"
$ws.Range("C56").Value = 3.9797
$ws.Range("D56").Value = 2.695
$ws.Range("E56").Value = 7.1162

# Row 57  # -> PVRINOX
$ws.Range("B57").Value = "PVRINOX"
$ws.Range("C57").Value = 3.9767
$ws.Range("D57").Value = 6.0724
$ws.Range("E57").Value = 14.5581

# Row 58  # -> SAMBHV
$ws.Range("B58").Value = "SAMBHV"
$ws.Range("C58").Value = 3.9713
$ws.Range("D58").Value = 2.4628
$ws.Range("E58").Value = 5.0018

# Row 59  # -> HITECHGEAR
$ws.Range("B59").Value = "HITECHGEAR"
$ws.Range("C59").Value = 3.9677
$ws.Range("D59").Value = 1.2548
$ws.Range("E59").Value = 10.0407

# Row 60  # -> JKIL
$ws.Range("B60").Value = "JKIL"
$ws.Range("C60").Value = 3.9042
$ws.Range("D60").Value = 2.716
$ws.Range("E60").Value = 1.5307

# Row 61  # -> GPIL
$ws.Range("B61").Value = "GPIL"
$ws.Range("C61").Value = 3.9001
$ws.Range("D61").Value = 6.065
$ws.Range("E61").Value = 14.1585

# Row 62  # -> HLEGLAS
$ws.Range("B62").Value = "HLEGLAS"
$ws.Range("C62").Value = 3.8678
$ws.Range("D62").Value = 8.333299999999999
$ws.Range("E62").Value = 27.38

# Row 63  # -> ORIENTTECH
$ws.Range("B63").Value = "ORIENTTECH"
$ws.Range("C63").Value = 3.8602
$ws.Range("D63").Value = 0.5569
$ws.Range("E63").Value = 32.7208

# Row 65  # -> INDIANHUME
$ws.Range("B65").Value = "INDIANHUME"
$ws.Range("C65").Value = 3.776
$ws.Range("D65").Value = 4.1198
$ws.Range("E65").Value = 0.6966

# Row 66  # -> SHK
$ws.Range("B66").Value = "SHK"
$ws.Range("C66").Value = 3.7171
$ws.Range("D66").Value = 2.4694
$ws.Range("E66").Value = -1.854

# Row 67
$ws.Range("C67").Value = 3.709
$ws.Range("D67").Value = 2.6111
$ws.Range("E67").Value = 7.3459

# Row 68
$ws.Range("C68").Value = 3.662
$ws.Range("D68").Value = 1.3094
$ws.Range("E68").Value = -1.2423

# Row 69  # -> GAIL
$ws.Range("B69").Value = "GAIL"
$ws.Range("C69").Value = 3.6311
$ws.Range("D69").Value = 2.1655
$ws.Range("E69").Value = 4.9067

# Row 70  # -> CMSINFO
$ws.Range("B70").Value = "CMSINFO"
$ws.Range("C70").Value = 3.6255
$ws.Range("D70").Value = 2.4064
$ws.Range("E70").Value = 2.6122

# Row 71  # -> MAITHANALL
$ws.Range("B71").Value = "MAITHANALL"
$ws.Range("C71").Value = 3.616
$ws.Range("D71").Value = 2.6571
$ws.Range("E71").Value = 1.9287

# Row 72  # -> RAJRATAN
$ws.Range("B72").Value = "RAJRATAN"
$ws.Range("C72").Value = 3.5627
$ws.Range("D72").Value = 1.109
$ws.Range("E72").Value = 27.1675

# Row 73  # -> ICRA
$ws.Range("B73").Value = "ICRA"
$ws.Range("C73").Value = 3.5362
$ws.Range("D73").Value = 4.2153
$ws.Range("E73").Value = 2.6229

# Row 74  # -> SALZERELEC
$ws.Range("B74").Value = "SALZERELEC"
$ws.Range("C74").Value = 3.5142
$ws.Range("D74").Value = 6.3309
$ws.Range("E74").Value = 16.7275

# Row 75  # -> AVALON
$ws.Range("B75").Value = "AVALON"
$ws.Range("C75").Value = 3.4721
$ws.Range("D75").Value = 8.236700000000001
$ws.Range("E75").Value = 20.1903

# Row 76  # -> PROSTARM
$ws.Range("B76").Value = "PROSTARM"
$ws.Range("C76").Value = 3.4327
$ws.Range("D76").Value = 0.5638
$ws.Range("E76").Value = -8.3277
